$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10-79 down to 11-80.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record.
$ws.Cells.Item(10, 1).Value = 11
$ws.Cells.Item(10, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(10, 3).Value = 'Bíobío'
$ws.Cells.Item(10, 4).Value = [datetime]'2023-04-19'
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 100112031
$ws.Cells.Item(10, 7).Value = 'Poroto verde'
$ws.Cells.Item(10, 8).Value = 'Magnum'
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 10).Value = 100
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 17000
$ws.Cells.Item(10, 13).Value = 16000
$ws.Cells.Item(10, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(10, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(10, 16).Value = 640
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = 'Hortaliza'

# Match the date display format used by the rest of the "Fecha" column.
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
